$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Move to location (9, 5) and remove the toolkit."
$ws.Range("B1").Value = "['Robot42']"
$ws.Range("E1").Value = "(9, 5)"

$ws.Range("A2").Value = "Move to location (3, 8) and remove the liquid spill."
$ws.Range("E2").Value = "(3, 8)"

$ws.Range("A3").Value = "Move to location (1, 4) and remove the large debris."
$ws.Range("B3").Value = "['Robot2', 'Robot39']"
$ws.Range("E3").Value = "(1, 4)"

$ws.Range("A4").Value = "Move to location (6, 5) and remove the dust."
$ws.Range("B4").Value = "['Robot50', 'Robot8']"
$ws.Range("E4").Value = "(6, 5)"

$ws.Range("A5").Value = "Move to location (9, 5) and remove the grass."
$ws.Range("B5").Value = "['Robot11']"
$ws.Range("E5").Value = "(9, 5)"

$ws.Range("A6").Value = "Move to location (5, 12) and remove the small debris."
$ws.Range("B6").Value = "['Robot28', 'Robot50']"
$ws.Range("E6").Value = "(5, 12)"

$ws.Range("A7").Value = "Move to location (11, 12) and remove the vehicle."
$ws.Range("E7").Value = "(11, 12)"

$ws.Range("A8").Value = "Move to location (12, 1) and remove the construction materials."
$ws.Range("B8").Value = "['Robot32', 'Robot29', 'Robot23']"
$ws.Range("E8").Value = "(12, 1)"

$ws.Range("A9").Value = "Move to location (8, 12) and remove the tree branches."
$ws.Range("E9").Value = "(8, 12)"

$ws.Range("A10").Value = "Move to location (3, 5) and remove the screws."
$ws.Range("E10").Value = "(3, 5)"

$wb.Save()
